# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates DAMSLTag (col I) and DialogAct (col J) values for the rows affected by the
# refreshed dialog-act annotation pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(8, 9).Value = "b"
$ws.Cells.Item(8, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(14, 9).Value = "b"
$ws.Cells.Item(14, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(19, 9).Value = "aa"
$ws.Cells.Item(19, 10).Value = "Agree/Accept"
$ws.Cells.Item(22, 9).Value = "sd"
$ws.Cells.Item(22, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(24, 9).Value = "sv"
$ws.Cells.Item(24, 10).Value = "Statement-opinion"
$ws.Cells.Item(30, 9).Value = "qy"
$ws.Cells.Item(30, 10).Value = "Yes-No-Question"
$ws.Cells.Item(39, 9).Value = "b"
$ws.Cells.Item(39, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(60, 9).Value = "aa"
$ws.Cells.Item(60, 10).Value = "Agree/Accept"
$ws.Cells.Item(83, 9).Value = "b"
$ws.Cells.Item(83, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(92, 9).Value = "b"
$ws.Cells.Item(92, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(113, 9).Value = "sd"
$ws.Cells.Item(113, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(126, 9).Value = "sv"
$ws.Cells.Item(126, 10).Value = "Statement-opinion"
$ws.Cells.Item(142, 9).Value = "sd"
$ws.Cells.Item(142, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(146, 9).Value = "aa"
$ws.Cells.Item(146, 10).Value = "Agree/Accept"
$ws.Cells.Item(157, 9).Value = "%"
$ws.Cells.Item(157, 10).Value = "Uninterpretable"
$ws.Cells.Item(165, 9).Value = "%"
$ws.Cells.Item(165, 10).Value = "Uninterpretable"
$ws.Cells.Item(173, 9).Value = "ba"
$ws.Cells.Item(173, 10).Value = "Appreciation"
$ws.Cells.Item(189, 9).Value = "ba"
$ws.Cells.Item(189, 10).Value = "Appreciation"
$ws.Cells.Item(207, 9).Value = "sd"
$ws.Cells.Item(207, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(223, 9).Value = "sd"
$ws.Cells.Item(223, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(225, 9).Value = "b"
$ws.Cells.Item(225, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(227, 9).Value = "aa"
$ws.Cells.Item(227, 10).Value = "Agree/Accept"
$ws.Cells.Item(228, 9).Value = "ba"
$ws.Cells.Item(228, 10).Value = "Appreciation"
$ws.Cells.Item(229, 9).Value = "b"
$ws.Cells.Item(229, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(259, 9).Value = "aa"
$ws.Cells.Item(259, 10).Value = "Agree/Accept"
$ws.Cells.Item(264, 9).Value = "sv"
$ws.Cells.Item(264, 10).Value = "Statement-opinion"
$ws.Cells.Item(283, 9).Value = "ba"
$ws.Cells.Item(283, 10).Value = "Appreciation"
$ws.Cells.Item(311, 9).Value = "sd"
$ws.Cells.Item(311, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(314, 9).Value = "aa"
$ws.Cells.Item(314, 10).Value = "Agree/Accept"
$ws.Cells.Item(323, 9).Value = "b"
$ws.Cells.Item(323, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(342, 9).Value = "sd"
$ws.Cells.Item(342, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(344, 9).Value = "sv"
$ws.Cells.Item(344, 10).Value = "Statement-opinion"
$ws.Cells.Item(345, 9).Value = "aa"
$ws.Cells.Item(345, 10).Value = "Agree/Accept"
$ws.Cells.Item(346, 9).Value = "b"
$ws.Cells.Item(346, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(348, 9).Value = "sd"
$ws.Cells.Item(348, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(351, 9).Value = "aa"
$ws.Cells.Item(351, 10).Value = "Agree/Accept"
$ws.Cells.Item(372, 9).Value = "sv"
$ws.Cells.Item(372, 10).Value = "Statement-opinion"
$ws.Cells.Item(384, 9).Value = "aa"
$ws.Cells.Item(384, 10).Value = "Agree/Accept"
$ws.Cells.Item(391, 9).Value = "b"
$ws.Cells.Item(391, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(397, 9).Value = "sv"
$ws.Cells.Item(397, 10).Value = "Statement-opinion"
$ws.Cells.Item(401, 9).Value = "sd"
$ws.Cells.Item(401, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(413, 9).Value = "sv"
$ws.Cells.Item(413, 10).Value = "Statement-opinion"
$ws.Cells.Item(414, 9).Value = "aa"
$ws.Cells.Item(414, 10).Value = "Agree/Accept"
$ws.Cells.Item(415, 9).Value = "sd"
$ws.Cells.Item(415, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(416, 9).Value = "sd"
$ws.Cells.Item(416, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(417, 9).Value = "sv"
$ws.Cells.Item(417, 10).Value = "Statement-opinion"
